$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Flip the RegressionTest answer (column D) from "Yes" to "No"
#        for rows 182 through 302. ---
$ws.Range("D182:D302").Value = "No"

# --- 2. Re-point the saved view: scroll position + active selection. ---
$excel.ActiveWindow.ScrollRow = 167
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C182").Select()

# --- 3. Rebuild the data-validation rules so D182:D302 carries its own
#        "Yes,No" list rule (separate from the D1:D181 block above it),
#        while E1 keeps its standalone "Yes"-only rule and E2:E302 / F2:F3
#        keep their "Yes,No" list rule. Re-creating E1 first and the
#        D/E/F "Yes,No" rules after matches the rule ordering produced
#        by re-validating D182:D302 in the live workbook. ---
$ws.Range("D1:D302").Validation.Delete()
$ws.Range("E2:E302").Validation.Delete()
$ws.Range("F2:F3").Validation.Delete()
$ws.Range("E1").Validation.Delete()

$ws.Range("E1").Validation.Add(3, 1, 3, """Yes""")
$ws.Range("D1:D181").Validation.Add(3, 1, 3, """Yes,No""")
$ws.Range("D182:D302").Validation.Add(3, 1, 3, """Yes,No""")
$ws.Range("E2:E302").Validation.Add(3, 1, 3, """Yes,No""")
$ws.Range("F2:F3").Validation.Add(3, 1, 3, """Yes,No""")
